# edit.ps1
# Applies the changes described by the diff:
#  1. Insert a new "Meta description" paragraph right after the title
#     (Heading1) paragraph, containing a bold "Meta description" run
#     followed by a plain run with the rest of the sentence.
#  2. Remove the trailing bold paragraph that duplicated the title text
#     ("Play Champion Raceway Free: Review & Gameplay Mechanics").
#  3. Replace the text of the trailing italic paragraph (which used to
#     hold the meta description) with the new feature-image prompt text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: insert the "Meta description" paragraph after the Heading1
# title paragraph, using InsertXML so the exact run layout (including
# the leading empty run) is reproduced faithfully.
# ---------------------------------------------------------------------
$nextPara = $d.Paragraphs.Item(2)

$targetRange = $d.Range($nextPara.Range.Start, $nextPara.Range.End)

$packagePrefix = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$packageSuffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$replacementXml = '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Discover the gameplay mechanics, pros, and cons of Champion Raceway in our review. Play for free and experience the unique multiplier system and two free spin features.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Gameplay Mechanics</w:t></w:r></w:p>'

$fullFragment = $packagePrefix + $replacementXml + $packageSuffix

$targetRange.InsertXML($fullFragment) | Out-Null

# ---------------------------------------------------------------------
# Step 2: remove the trailing duplicate-title paragraph (bold run) and
# replace the text of the italic paragraph that follows it.
# ---------------------------------------------------------------------
$dupTitlePara = $null
$imgPromptPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`n", "`x07")
    $styleName = $p.Style.NameLocal
    if ($t -eq "Play Champion Raceway Free: Review & Gameplay Mechanics" -and $styleName -ne "Heading 1") {
        $dupTitlePara = $p
    }
    if ($t -eq "Discover the gameplay mechanics, pros, and cons of Champion Raceway in our review. Play for free and experience the unique multiplier system and two free spin features.") {
        $imgPromptPara = $p
    }
}

# The italic paragraph's text is replaced with the new image-prompt text,
# preserving its existing run/paragraph formatting (italic). This is done
# BEFORE deleting the duplicate-title paragraph below it, so that the
# Range positions used here are not invalidated by the earlier deletion.
$imgPromptPara.Range.Find.Execute(
    "Discover the gameplay mechanics, pros, and cons of Champion Raceway in our review. Play for free and experience the unique multiplier system and two free spin features.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Create a cartoon-style feature image for Champion Raceway that features a happy Maya warrior with glasses. The image should showcase the horse racing theme of the game, with the Maya warrior standing next to a horse on a racetrack. Use bright colors to make the image eye-catching and visually appealing. Add in elements from the game, such as the Wild symbol or the racetrack above the reels, to tie it back to the game. Overall, the feature image should capture the fun and excitement of horse racing while also highlighting the unique aspects of Champion Raceway.",
    2) | Out-Null

# The duplicate title paragraph (bold run) is removed entirely.
$dupTitlePara.Range.Delete()
